# Work against the active workbook (the instructions' $wb placeholder is not
# populated by this runtime -- use $excel.ActiveWorkbook instead).
$wb = $excel.ActiveWorkbook

# --- 1. Update the "time_taken" timestamps on the "data" sheet (F2:F25) ---
$wsData = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:23:07.559507",
    "2021-10-05 14:23:07.559516",
    "2021-10-05 14:23:07.559519",
    "2021-10-05 14:23:07.559522",
    "2021-10-05 14:23:07.559525",
    "2021-10-05 14:23:07.559528",
    "2021-10-05 14:23:07.559530",
    "2021-10-05 14:23:07.559533",
    "2021-10-05 14:23:07.559536",
    "2021-10-05 14:23:07.559539",
    "2021-10-05 14:23:07.559542",
    "2021-10-05 14:23:07.559545",
    "2021-10-05 14:23:07.559547",
    "2021-10-05 14:23:07.559550",
    "2021-10-05 14:23:07.559553",
    "2021-10-05 14:23:07.559556",
    "2021-10-05 14:23:07.559559",
    "2021-10-05 14:23:07.559562",
    "2021-10-05 14:23:07.559564",
    "2021-10-05 14:23:07.559567",
    "2021-10-05 14:23:07.559570",
    "2021-10-05 14:23:07.559573",
    "2021-10-05 14:23:07.559576",
    "2021-10-05 14:23:07.559578"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add a new "metadata" worksheet placed after "data" ---
$wsMeta = $wb.Worksheets.Add($null, $wsData)
$wsMeta.Name = "metadata"

# Header row
$wsMeta.Cells.Item(1, 2).Value = "data_name"
$wsMeta.Cells.Item(1, 3).Value = "data_id"
$wsMeta.Cells.Item(1, 4).Value = "data_version"
$wsMeta.Cells.Item(1, 5).Value = "data_version_created"
$wsMeta.Cells.Item(1, 6).Value = "panel_query_time"
$wsMeta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$wsMeta.Cells.Item(2, 1).Value = 0
$wsMeta.Cells.Item(2, 2).Value = "Viral resistance"
$wsMeta.Cells.Item(2, 3).Value = 928
# "0.63" must stay a text string (not be coerced to a number) - force text format first.
$wsMeta.Range("D2").NumberFormat = "@"
$wsMeta.Cells.Item(2, 4).Value = "0.63"
$wsMeta.Cells.Item(2, 5).Value = "2020-07-07T11:06:18.864817Z"
$wsMeta.Cells.Item(2, 6).Value = "2021-10-05 14:23:07.555921"
$wsMeta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/928/?format=json"

# Copy the bold/bordered header style from the "data" sheet's header row,
# and the index-column style used for column A.
$wsData.Range("B1:F1").Copy()
$wsMeta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$wsData.Range("A2").Copy()
$wsMeta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsData.Activate()
